# Auto-applied updates to currentAveragePrice / Leve profit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 16680165
$ws.Range("J97").Value = 16680165
$ws.Range("L97").Value = 50040495
$ws.Range("N97").Value = -50041487

$ws.Range("H106").Value = 1747.2
$ws.Range("I106").Value = 1747.2
$ws.Range("K106").Value = 1747.2
$ws.Range("M106").Value = -1116.2

$ws.Range("H112").Value = 9392.305
$ws.Range("J112").Value = 10065.143
$ws.Range("L112").Value = 30195.429
$ws.Range("N112").Value = -32411.429

$ws.Range("H132").Value = 1345.8776
$ws.Range("I132").Value = 1385.1522
$ws.Range("J132").Value = 743.6667
$ws.Range("K132").Value = 4155.4566
$ws.Range("L132").Value = 2231.0001
$ws.Range("M132").Value = -1625.4566
$ws.Range("N132").Value = -7291.0001

$ws.Range("H138").Value = 3338674.8
$ws.Range("I138").Value = 2152.7273
$ws.Range("J138").Value = 5270345
$ws.Range("K138").Value = 6458.1819
$ws.Range("L138").Value = 15811035
$ws.Range("M138").Value = -1318.1819
$ws.Range("N138").Value = -15821315

$ws.Range("H141").Value = 71430320
$ws.Range("I141").Value = 83334580
$ws.Range("K141").Value = 250003740
$ws.Range("M141").Value = -249998560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2980517.8
$ws.Range("I32").Value = 3381726.2
$ws.Range("J32").Value = 11576.2
$ws.Range("K32").Value = 3381726.2
$ws.Range("L32").Value = 11576.2
$ws.Range("M32").Value = -3381439.2
$ws.Range("N32").Value = -12150.2

$ws.Range("H61").Value = 6347.5137
$ws.Range("I61").Value = 2353.1667
$ws.Range("J61").Value = 13721.692
$ws.Range("K61").Value = 2353.1667
$ws.Range("L61").Value = 13721.692
$ws.Range("M61").Value = -2141.1667
$ws.Range("N61").Value = -14145.692

$ws.Range("H74").Value = 38249.17
$ws.Range("I74").Value = 57946.11
$ws.Range("J74").Value = 6017.8184
$ws.Range("K74").Value = 57946.11
$ws.Range("L74").Value = 6017.8184
$ws.Range("M74").Value = -57072.11
$ws.Range("N74").Value = -7765.8184

$ws.Range("H77").Value = 38249.17
$ws.Range("I77").Value = 57946.11
$ws.Range("J77").Value = 6017.8184
$ws.Range("K77").Value = 289730.55
$ws.Range("L77").Value = 30089.092
$ws.Range("M77").Value = -285362.55
$ws.Range("N77").Value = -38825.092

$ws.Range("H97").Value = 9276731
$ws.Range("I97").Value = 791.5
$ws.Range("J97").Value = 11927000
$ws.Range("K97").Value = 791.5
$ws.Range("L97").Value = 11927000
$ws.Range("M97").Value = -295.5
$ws.Range("N97").Value = -11927992

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws.Range("H132").Value = 1256266.8
$ws.Range("I132").Value = 3578668
$ws.Range("J132").Value = 5743.0386
$ws.Range("K132").Value = 10736004
$ws.Range("L132").Value = 17229.1158
$ws.Range("M132").Value = -10733474
$ws.Range("N132").Value = -22289.1158

$ws.Range("H136").Value = 6347.5137
$ws.Range("I136").Value = 2353.1667
$ws.Range("J136").Value = 13721.692
$ws.Range("K136").Value = 7059.500100000001
$ws.Range("L136").Value = 41165.076
$ws.Range("M136").Value = -4509.500100000001
$ws.Range("N136").Value = -46265.076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 40003380
$ws.Range("I94").Value = 62501340
$ws.Range("K94").Value = 62501340
$ws.Range("M94").Value = -62500889

$ws.Range("H103").Value = 26666
$ws.Range("J103").Value = 26666
$ws.Range("L103").Value = 26666
$ws.Range("N103").Value = -29010

$ws.Range("H124").Value = 51549
$ws.Range("J124").Value = 51549
$ws.Range("L124").Value = 51549
$ws.Range("N124").Value = -61369

$ws.Range("H126").Value = 42145
$ws.Range("J126").Value = 51590
$ws.Range("L126").Value = 51590
$ws.Range("N126").Value = -61470

$ws.Range("H134").Value = 6055.129
$ws.Range("I134").Value = 2075.2354
$ws.Range("J134").Value = 10887.857
$ws.Range("K134").Value = 6225.706200000001
$ws.Range("L134").Value = 32663.571
$ws.Range("M134").Value = -3690.706200000001
$ws.Range("N134").Value = -37733.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7702.8
$ws.Range("I31").Value = 2040
$ws.Range("J31").Value = 8835.360000000001
$ws.Range("K31").Value = 2040
$ws.Range("L31").Value = 8835.360000000001
$ws.Range("M31").Value = -1745
$ws.Range("N31").Value = -9425.360000000001

$ws.Range("H34").Value = 7702.8
$ws.Range("I34").Value = 2040
$ws.Range("J34").Value = 8835.360000000001
$ws.Range("K34").Value = 2040
$ws.Range("L34").Value = 8835.360000000001
$ws.Range("M34").Value = -1838
$ws.Range("N34").Value = -9239.360000000001

$ws.Range("H58").Value = 7585.75
$ws.Range("I58").Value = 1861.875
$ws.Range("J58").Value = 9875.299999999999
$ws.Range("K58").Value = 1861.875
$ws.Range("L58").Value = 9875.299999999999
$ws.Range("M58").Value = -1658.875
$ws.Range("N58").Value = -10281.3

$ws.Range("H135").Value = 49994
$ws.Range("J135").Value = 49994
$ws.Range("L135").Value = 49994
$ws.Range("N135").Value = -60134

$ws.Range("H136").Value = 7585.75
$ws.Range("I136").Value = 1861.875
$ws.Range("J136").Value = 9875.299999999999
$ws.Range("K136").Value = 5585.625
$ws.Range("L136").Value = 29625.9
$ws.Range("M136").Value = -3035.625
$ws.Range("N136").Value = -34725.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2500.6667
$ws.Range("I98").Value = 1626.25
$ws.Range("K98").Value = 4878.75
$ws.Range("M98").Value = -3380.75

$ws.Range("H129").Value = 8384116.5
$ws.Range("I129").Value = 444.0909
$ws.Range("J129").Value = 18630828
$ws.Range("K129").Value = 1332.2727
$ws.Range("L129").Value = 55892484
$ws.Range("M129").Value = 3667.7273
$ws.Range("N129").Value = -55902484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2110.7144
$ws.Range("I97").Value = 2263.7778
$ws.Range("K97").Value = 2263.7778
$ws.Range("M97").Value = -1767.7778

$ws.Range("H102").Value = 3263.8647
$ws.Range("I102").Value = 3289.0857
$ws.Range("K102").Value = 3289.0857
$ws.Range("M102").Value = -1667.0857

$ws.Range("H132").Value = 4191.1
$ws.Range("I132").Value = 2293.0667
$ws.Range("J132").Value = 6089.1333
$ws.Range("K132").Value = 6879.2001
$ws.Range("L132").Value = 18267.3999
$ws.Range("M132").Value = -4349.2001
$ws.Range("N132").Value = -23327.3999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7193.9473
$ws.Range("I7").Value = 6337.375
$ws.Range("J7").Value = 7816.909
$ws.Range("K7").Value = 6337.375
$ws.Range("L7").Value = 7816.909
$ws.Range("M7").Value = -6225.375
$ws.Range("N7").Value = -8040.909

$ws.Range("H22").Value = 3893.8572
$ws.Range("I22").Value = 1371.5
$ws.Range("K22").Value = 1371.5
$ws.Range("M22").Value = -1076.5

$ws.Range("H27").Value = 3893.8572
$ws.Range("I27").Value = 1371.5
$ws.Range("K27").Value = 1371.5
$ws.Range("M27").Value = -1264.5

$ws.Range("H46").Value = 3477.5454
$ws.Range("I46").Value = 1937.75
$ws.Range("J46").Value = 4357.4287
$ws.Range("K46").Value = 1937.75
$ws.Range("L46").Value = 4357.4287
$ws.Range("M46").Value = -1749.75
$ws.Range("N46").Value = -4733.4287

$ws.Range("H126").Value = 7193.9473
$ws.Range("I126").Value = 6337.375
$ws.Range("J126").Value = 7816.909
$ws.Range("K126").Value = 19012.125
$ws.Range("L126").Value = 23450.727
$ws.Range("M126").Value = -16542.125
$ws.Range("N126").Value = -28390.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8404305
$ws.Range("I81").Value = 589566
$ws.Range("J81").Value = 25010624
$ws.Range("K81").Value = 1179132
$ws.Range("L81").Value = 50021248
$ws.Range("M81").Value = -1178071
$ws.Range("N81").Value = -50023370

$ws.Range("H84").Value = 8404305
$ws.Range("I84").Value = 589566
$ws.Range("J84").Value = 25010624
$ws.Range("K84").Value = 5895660
$ws.Range("L84").Value = 250106240
$ws.Range("M84").Value = -5890356
$ws.Range("N84").Value = -250116848

$ws.Range("H96").Value = 4524.25
$ws.Range("I96").Value = 4199
$ws.Range("J96").Value = 5500
$ws.Range("K96").Value = 4199
$ws.Range("L96").Value = 5500
$ws.Range("M96").Value = -2826
$ws.Range("N96").Value = -8246

$ws.Range("H122").Value = 142593.77
$ws.Range("I122").Value = 164048.8
$ws.Range("J122").Value = 8499.75
$ws.Range("K122").Value = 492146.4
$ws.Range("L122").Value = 25499.25
$ws.Range("M122").Value = -489696.4
$ws.Range("N122").Value = -30399.25
